# "include update roommate feature"
# Update the Sprint Stories tracker:
#   - "I want to view list of current room"  -> Status: IN PROGRESS -> DONE
#   - "I want to update room detail"         -> Status: NOT STARTED -> DONE
#   - "I want to add new expense for this month" -> Status: NOT STARTED -> IN PROGRESS
# and move the active selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C5 needs to become "IN PROGRESS", styled the same way C3 currently is
# (the "Neutral" table style with border used for IN PROGRESS). Copy that
# formatting over to C5 first, before C3's own formatting is overwritten.
$ws.Range("C3").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "IN PROGRESS"

# C3 and C4 both need to become "DONE", styled the same way C2 currently is
# (the "Good" table style with border used for DONE).
$ws.Range("C2").Copy()

$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "DONE"

$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "DONE"

# Move / record the active selection on C5.
$ws.Range("C5").Select()
